$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new product row (id 7) to the table
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "vicio"
$ws.Range("C8").Value = "imagen de producto"
$ws.Range("D8").Value = 1.3
$ws.Range("E8").Value = "S/"

# Update selection to match the post-edit state
$ws.Range("B15").Select()
